# Update ExampleWQX.xlsx: add "Template updated" and "Samples updated" notes
# (new red text in column C of the Instructions sheet) and move the active
# selection on both sheets, matching the upstream "update all templates and
# sample files in inst #48" commit.

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Meta")
$wsInstructions = $wb.Worksheets.Item("Instructions")

# New dated annotations, shown in red text next to the existing notes on the
# Instructions tab.
$wsInstructions.Range("C1").Value = "Template updated 1/17/23"
$wsInstructions.Range("C1").Font.Color = 255

$wsInstructions.Range("C2").Value = "Samples updated 1/8/23"
$wsInstructions.Range("C2").Font.Color = 255

# Update the remembered selection on each sheet, then leave the Meta sheet
# active/selected as it was before.
[void]$wsInstructions.Range("C3").Select()
[void]$wsMeta.Range("B12").Select()
[void]$wsMeta.Activate()
